# Build site at 2022-01-09 00:29:46 UTC
# Insert two new rows describing the instructors ("Docentes responsáveis:")
# right after the "Objectives:" block (currently rows 10-11) and before the
# "Programa resumido:" block, shifting everything below down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 12, 13, 14 (existing rows 12-20 shift to 15-23).
$ws.Rows("12:14").Insert()

# Row 12: section header only (column A), matching the style of the other
# single-column header rows (e.g. row 19 "Avaliação:").
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B12:C12").Clear()

# Row 13: first instructor, duplicated into B and C (like every other
# content row on this sheet), no value in column A.
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C13").Value = "198273 - Domingos Savio Giordani"

# Row 14: second instructor, same pattern.
$ws.Range("A14").Clear()
$ws.Range("B14").Value = "1506103 - Pedro Carlos de Oliveira"
$ws.Range("C14").Value = "1506103 - Pedro Carlos de Oliveira"
